$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 623.8
$ws.Range("I12").Value = 769
$ws.Range("J12").Value = 527
$ws.Range("K12").Value = 769
$ws.Range("L12").Value = 527
$ws.Range("M12").Value = -599
$ws.Range("N12").Value = -867
$ws.Range("H17").Value = 843.1818
$ws.Range("J17").Value = 870.8095
$ws.Range("L17").Value = 2612.4285
$ws.Range("N17").Value = -2948.4285
$ws.Range("H33").Value = 172.77777
$ws.Range("I33").Value = 175.625
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 175.625
$ws.Range("L33").Value = 150
$ws.Range("M33").Value = 53.375
$ws.Range("N33").Value = -608
$ws.Range("H43").Value = 1193.3334
$ws.Range("I43").Value = 1190
$ws.Range("K43").Value = 1190
$ws.Range("M43").Value = -1121
$ws.Range("H53").Value = 350.45456
$ws.Range("I53").Value = 259.83334
$ws.Range("K53").Value = 259.83334
$ws.Range("M53").Value = 377.16666
$ws.Range("H70").Value = 1661.875
$ws.Range("I70").Value = 1397.5
$ws.Range("K70").Value = 4192.5
$ws.Range("M70").Value = -3922.5
$ws.Range("H73").Value = 1661.875
$ws.Range("I73").Value = 1397.5
$ws.Range("K73").Value = 4192.5
$ws.Range("M73").Value = -3256.5
$ws.Range("H82").Value = 241
$ws.Range("I82").Value = 241
$ws.Range("K82").Value = 723
$ws.Range("M82").Value = -317
$ws.Range("H85").Value = 241
$ws.Range("I85").Value = 241
$ws.Range("K85").Value = 723
$ws.Range("M85").Value = 681
$ws.Range("H98").Value = 974.375
$ws.Range("I98").Value = 1074.25
$ws.Range("J98").Value = 874.5
$ws.Range("K98").Value = 1074.25
$ws.Range("L98").Value = 874.5
$ws.Range("M98").Value = 423.75
$ws.Range("N98").Value = -3870.5
$ws.Range("H103").Value = 883.3333
$ws.Range("J103").Value = 1050
$ws.Range("L103").Value = 3150
$ws.Range("N103").Value = -4322
$ws.Range("H122").Value = 974.375
$ws.Range("I122").Value = 1074.25
$ws.Range("J122").Value = 874.5
$ws.Range("K122").Value = 3222.75
$ws.Range("L122").Value = 2623.5
$ws.Range("M122").Value = -772.75
$ws.Range("N122").Value = -7523.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 470.16666
$ws.Range("I5").Value = 470.16666
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 470.16666
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -358.16666
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 5927565.5
$ws.Range("I32").Value = 5838196
$ws.Range("J32").Value = 7000000
$ws.Range("K32").Value = 5838196
$ws.Range("L32").Value = 7000000
$ws.Range("M32").Value = -5837909
$ws.Range("N32").Value = -7000574
$ws.Range("H74").Value = 2799.8572
$ws.Range("I74").Value = 2799.8572
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2799.8572
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1925.8572
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2799.8572
$ws.Range("I77").Value = 2799.8572
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 13999.286
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -9631.286
$ws.Range("N77").ClearContents()
$ws.Range("H97").Value = 677.8182
$ws.Range("J97").Value = 673.25
$ws.Range("L97").Value = 673.25
$ws.Range("N97").Value = -1665.25
$ws.Range("H122").Value = 1033.3334
$ws.Range("I122").Value = 1033.3334
$ws.Range("K122").Value = 3100.0002
$ws.Range("M122").Value = -650.0001999999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 470.16666
$ws.Range("I4").Value = 470.16666
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 470.16666
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -355.16666
$ws.Range("N4").ClearContents()
$ws.Range("H20").Value = 2300.3635
$ws.Range("I20").Value = 1651.8334
$ws.Range("J20").Value = 3078.6
$ws.Range("K20").Value = 1651.8334
$ws.Range("L20").Value = 3078.6
$ws.Range("M20").Value = -1404.8334
$ws.Range("N20").Value = -3572.6
$ws.Range("H134").Value = 2158.5833
$ws.Range("I134").Value = 2158.5833
$ws.Range("K134").Value = 6475.749899999999
$ws.Range("M134").Value = -3940.749899999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3925.5334
$ws.Range("I16").Value = 2567.375
$ws.Range("K16").Value = 2567.375
$ws.Range("M16").Value = -2280.375
$ws.Range("H31").Value = 2943.25
$ws.Range("I31").Value = 2444
$ws.Range("K31").Value = 2444
$ws.Range("M31").Value = -2149
$ws.Range("H34").Value = 2943.25
$ws.Range("I34").Value = 2444
$ws.Range("K34").Value = 2444
$ws.Range("M34").Value = -2242
$ws.Range("H58").Value = 2333.25
$ws.Range("H113").Value = 3925.5334
$ws.Range("I113").Value = 2567.375
$ws.Range("K113").Value = 2567.375
$ws.Range("M113").Value = -397.375
$ws.Range("H136").Value = 2333.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1516.5
$ws.Range("J45").Value = 1516.5
$ws.Range("L45").Value = 4549.5
$ws.Range("N45").Value = -5613.5
$ws.Range("H56").Value = 11021.155
$ws.Range("I56").Value = 11021.155
$ws.Range("K56").Value = 11021.155
$ws.Range("M56").Value = -10491.155
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4881.2
$ws.Range("I70").Value = 4599.5
$ws.Range("K70").Value = 4599.5
$ws.Range("M70").Value = -4329.5
$ws.Range("H73").Value = 4881.2
$ws.Range("I73").Value = 4599.5
$ws.Range("K73").Value = 4599.5
$ws.Range("M73").Value = -3663.5
$ws.Range("H80").Value = 4966.5
$ws.Range("I80").Value = 2166.3333
$ws.Range("K80").Value = 2166.3333
$ws.Range("M80").Value = -1168.3333
$ws.Range("H83").Value = 4966.5
$ws.Range("I83").Value = 2166.3333
$ws.Range("K83").Value = 10831.6665
$ws.Range("M83").Value = -5839.666499999999
$ws.Range("H94").Value = 31159.5
$ws.Range("J94").Value = 31159.5
$ws.Range("L94").Value = 31159.5
$ws.Range("N94").Value = -32511.5
$ws.Range("H102").Value = 2276.5386
$ws.Range("I102").Value = 2281.0833
$ws.Range("K102").Value = 2281.0833
$ws.Range("M102").Value = -659.0832999999998
$ws.Range("H113").Value = 1045
$ws.Range("I113").Value = 799
$ws.Range("K113").Value = 799
$ws.Range("M113").Value = 1371

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 14991.333
$ws.Range("H20").Value = 11866.167
$ws.Range("I20").Value = 8733
$ws.Range("K20").Value = 8733
$ws.Range("M20").Value = -8507
$ws.Range("H22").Value = 871.44446
$ws.Range("I22").Value = 849.7143
$ws.Range("K22").Value = 849.7143
$ws.Range("M22").Value = -554.7143
$ws.Range("H27").Value = 871.44446
$ws.Range("I27").Value = 849.7143
$ws.Range("K27").Value = 849.7143
$ws.Range("M27").Value = -742.7143
$ws.Range("H38").Value = 5000
$ws.Range("I38").Value = 5000
$ws.Range("K38").Value = 5000
$ws.Range("M38").Value = -4590
$ws.Range("H39").Value = 14999
$ws.Range("J39").Value = 14999
$ws.Range("L39").Value = 14999
$ws.Range("N39").Value = -15919
$ws.Range("H46").Value = 2533.2
$ws.Range("I46").Value = 1666.5
$ws.Range("K46").Value = 1666.5
$ws.Range("M46").Value = -1478.5
$ws.Range("H68").Value = 860.1429000000001
$ws.Range("I68").Value = 841
$ws.Range("K68").Value = 841
$ws.Range("M68").Value = -92
$ws.Range("H71").Value = 860.1429000000001
$ws.Range("I71").Value = 841
$ws.Range("K71").Value = 4205
$ws.Range("M71").Value = -461

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5736.4443
$ws.Range("I14").Value = 732.75
$ws.Range("J14").Value = 9739.4
$ws.Range("K14").Value = 732.75
$ws.Range("L14").Value = 9739.4
$ws.Range("M14").Value = -564.75
$ws.Range("N14").Value = -10075.4
$ws.Range("H27").Value = 39990
$ws.Range("J27").Value = 39990
$ws.Range("L27").Value = 39990
$ws.Range("N27").Value = -40128
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H30").Value = 4553.8
$ws.Range("I30").Value = 3999.5
$ws.Range("K30").Value = 3999.5
$ws.Range("M30").Value = -3892.5
$ws.Range("H62").Value = 20400
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 27333.334
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 27333.334
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -28581.334
$ws.Range("H65").Value = 20400
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 27333.334
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 136666.67
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -142906.67
$ws.Range("H81").Value = 5002500.5
$ws.Range("J81").Value = 5002500.5
$ws.Range("L81").Value = 10005001
$ws.Range("N81").Value = -10007123
$ws.Range("H84").Value = 5002500.5
$ws.Range("J84").Value = 5002500.5
$ws.Range("L84").Value = 50025005
$ws.Range("N84").Value = -50035613
